# Applies the diff: adds column I ("Macro Weather Only") and updates
# the recomputed B:H statistic values for the existing strategy columns.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing B:H values that were recomputed ---
$ws.Range("B6").Value = 97952.11323940507
$ws.Range("C6").Value = 54416.69230427477
$ws.Range("D6").Value = 78790.7765345711
$ws.Range("E6").Value = 61372.32777811344
$ws.Range("F6").Value = 95447.31550612622
$ws.Range("G6").Value = 95254.50356179113
$ws.Range("H6").Value = 98705.10827782101
$ws.Range("B7").Value = 879.5211323940507
$ws.Range("C7").Value = 444.1669230427476
$ws.Range("D7").Value = 687.9077653457109
$ws.Range("E7").Value = 513.7232777811345
$ws.Range("F7").Value = 854.4731550612622
$ws.Range("G7").Value = 852.5450356179114
$ws.Range("H7").Value = 887.0510827782102
$ws.Range("B10").Value = 291.2621359223303
$ws.Range("C10").Value = 45362.14682513379
$ws.Range("D10").Value = 39802.22156262527
$ws.Range("E10").Value = 31292.23019318161
$ws.Range("F10").Value = 32367.25999974222
$ws.Range("G10").Value = 34533.47606432099
$ws.Range("H10").Value = 39741.31627493749
$ws.Range("B11").Value = 76.67018162770141
$ws.Range("C11").Value = 69.85770800230848
$ws.Range("D11").Value = 60.03165909794801
$ws.Range("E11").Value = 64.98308446010253
$ws.Range("F11").Value = 50.00665535685685
$ws.Range("G11").Value = 48.59852711443431
$ws.Range("H11").Value = 41.02262800007097
$ws.Range("C12").Value = 1423
$ws.Range("D12").Value = 1062
$ws.Range("E12").Value = 1423
$ws.Range("F12").Value = 1048
$ws.Range("G12").Value = 1048
$ws.Range("H12").Value = 963
$ws.Range("B16").Value = 87952.11323940507
$ws.Range("E16").Value = 6873.12386019776
$ws.Range("C17").Value = 21.05263157894737
$ws.Range("D17").Value = 28.57142857142857
$ws.Range("E17").Value = 42.30769230769231
$ws.Range("F17").Value = 28.57142857142857
$ws.Range("G17").Value = 38.88888888888889
$ws.Range("H17").Value = 50
$ws.Range("C18").Value = 580.3436609847811
$ws.Range("D18").Value = 577.9970618111165
$ws.Range("E18").Value = 579.2055603855536
$ws.Range("F18").Value = 577.9970618111165
$ws.Range("G18").Value = 577.9970618111165
$ws.Range("H18").Value = 577.9970618111165
$ws.Range("C19").Value = -16.62711095565251
$ws.Range("D19").Value = -16.6271109556525
$ws.Range("E19").Value = -33.68688198368736
$ws.Range("F19").Value = -16.62711095565251
$ws.Range("G19").Value = -16.62711095565251
$ws.Range("H19").Value = -16.62711095565252
$ws.Range("C20").Value = 183.2799761902225
$ws.Range("D20").Value = 180.6338657263656
$ws.Range("E20").Value = 110.142744204652
$ws.Range("F20").Value = 180.6338657263657
$ws.Range("G20").Value = 115.7343873924318
$ws.Range("H20").Value = 70.57729564203204
$ws.Range("C21").Value = -8.508381623259538
$ws.Range("D21").Value = -8.67210319871198
$ws.Range("E21").Value = -9.954409587018603
$ws.Range("F21").Value = -8.672103198711978
$ws.Range("G21").Value = -9.56799671808717
$ws.Range("H21").Value = -10.41293652185093
$ws.Range("C22").Value = 241.5
$ws.Range("D22").Value = 240.25
$ws.Range("E22").Value = 247.7272727272685
$ws.Range("F22").Value = 240.25
$ws.Range("G22").Value = 164.285714285706
$ws.Range("H22").Value = 137.0833333333333
$ws.Range("C23").Value = 23.2
$ws.Range("D23").Value = 28
$ws.Range("E23").Value = 15.73333333333333
$ws.Range("F23").Value = 28
$ws.Range("G23").Value = 25.54545454545139
$ws.Range("H23").Value = 23.91666666666667
$ws.Range("C24").Value = 1.820471186574542
$ws.Range("D24").Value = 2.526261394793529
$ws.Range("E24").Value = 2.493161574345239
$ws.Range("F24").Value = 3.375617835108335
$ws.Range("G24").Value = 3.480337368829281
$ws.Range("H24").Value = 3.8818334294992
$ws.Range("C25").Value = 2337.720647593409
$ws.Range("D25").Value = 4913.626895326507
$ws.Range("E25").Value = 1711.507842996757
$ws.Range("F25").Value = 6103.379679009016
$ws.Range("G25").Value = 4736.361308988395
$ws.Range("H25").Value = 3696.046178242543
$ws.Range("B26").Value = 1.032371909317504
$ws.Range("C26").Value = 0.9166960627269037
$ws.Range("D26").Value = 1.074490484597722
$ws.Range("E26").Value = 0.960067769993117
$ws.Range("F26").Value = 1.224484835567919
$ws.Range("G26").Value = 1.277343028236645
$ws.Range("H26").Value = 1.333544678406093
$ws.Range("B27").Value = 0.7512559824029422
$ws.Range("C27").Value = 0.5750683938804114
$ws.Range("D27").Value = 0.8479918939357475
$ws.Range("E27").Value = 0.6705554229278602
$ws.Range("F27").Value = 1.135593004904166
$ws.Range("G27").Value = 1.16719627245031
$ws.Range("H27").Value = 1.409946318730348
$ws.Range("B28").Value = 1.175629984046056
$ws.Range("C28").Value = 1.173787098595113
$ws.Range("D28").Value = 1.213956251052636
$ws.Range("E28").Value = 1.180455168000282
$ws.Range("F28").Value = 1.260055828077987
$ws.Range("G28").Value = 1.281919577265156
$ws.Range("H28").Value = 1.317592467188496
$ws.Range("B29").Value = 1.506382083567063
$ws.Range("C29").Value = 1.425907363880873
$ws.Range("D29").Value = 1.68432769177371
$ws.Range("E29").Value = 1.476718805724387
$ws.Range("F29").Value = 1.971574292753851
$ws.Range("G29").Value = 2.085044207408006
$ws.Range("H29").Value = 2.201453872372344

# --- Step 2: add new column I ("Macro Weather Only") ---

# Header: copy style from H1 (bold, bordered, centered) onto I1, then set text
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "Macro Weather Only"

# I2 / I3: copy date-format style from H2 / H3
$ws.Range("H2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$ws.Range("I2").Value = 43891
$ws.Range("H3").Copy() | Out-Null
$ws.Range("I3").PasteSpecial(-4122) | Out-Null
$ws.Range("I3").Value = 45721

# I4: copy integer-format style from H4
$ws.Range("H4").Copy() | Out-Null
$ws.Range("I4").PasteSpecial(-4122) | Out-Null
$ws.Range("I4").Value = 1831

# I5:I11 - plain values, no special number format
$ws.Range("I5").Value = 10000
$ws.Range("I6").Value = 110801.8459452796
$ws.Range("I7").Value = 1008.018459452796
$ws.Range("I8").Value = 908.9067663658728
$ws.Range("I9").Value = 100
$ws.Range("I10").Value = 32252.11477624038
$ws.Range("I11").Value = 53.85253498697061

# I12: copy integer-format style from H12
$ws.Range("H12").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("I12").Value = 702

# I13:I21 - plain values
$ws.Range("I13").Value = 13
$ws.Range("I14").Value = 12
$ws.Range("I15").Value = 1
$ws.Range("I16").Value = -3191.791299165622
$ws.Range("I17").Value = 75
$ws.Range("I18").Value = 393.0063986615069
$ws.Range("I19").Value = -8.195447894889046
$ws.Range("I20").Value = 98.01933666509163
$ws.Range("I21").Value = -5.450345487851835

# I22 / I23: copy integer-format style from H22 / H23
$ws.Range("H22").Copy() | Out-Null
$ws.Range("I22").PasteSpecial(-4122) | Out-Null
$ws.Range("I22").Value = 193.8888888888889
$ws.Range("H23").Copy() | Out-Null
$ws.Range("I23").PasteSpecial(-4122) | Out-Null
$ws.Range("I23").Value = 60.66666666666666

# I24:I29 - plain values
$ws.Range("I24").Value = 11.68137008168624
$ws.Range("I25").Value = 8666.136437037103
$ws.Range("I26").Value = 1.21437771306743
$ws.Range("I27").Value = 1.142368674134976
$ws.Range("I28").Value = 1.290232934790686
$ws.Range("I29").Value = 1.807680895864954

# Clear clipboard/marching-ants marquee left over from the copy operations
$excel.CutCopyMode = $false

Write-Output "Done: column I added, B:H stats refreshed."